$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A35").Value = 45881
$ws.Range("A35").NumberFormat = $ws.Range("A34").NumberFormat

$ws.Range("B35").Value = 6441.64013671875
$ws.Range("C35").Value = 6444.25
$ws.Range("D35").Value = 6385.759765625
$ws.Range("E35").Value = 6395.169921875
$ws.Range("F35").Value = 1667322000
$ws.Range("G35").Value = 0.0106990624099332
$ws.Range("H35").Value = 1

$ws.Range("I35").Value = @'
from power plant to the moon the u nuclear revival in dogecoin jump in a week what is going on teenage mutant ninja turtle pizzeria to debut in la thestreet pro analyst tackle tariff crypto and more in market roundtable solar stock outshine sp in despite policy headwind why share of tesla are soaring today time tracking and trust what today employee really think amazon is selling a storage cart for thats great for small place cold wallet m presale entry beat xrps price test and eths resistance goal why many smallbusinesses owner fear even a short recession will sink them creating a load history tracker that improves rate negotiation crypto news today strategy m bitcoin purchase bullish upsizes ipo alt sigma b fundraising move solana price prediction etf approval could help sol overtake xrp but this cheaper token will rally x and outshine their roi bank of america strategist expect fed to lower interest rate benny brown acquires wide open music publishing form new publishing company affirm stock revisits this place a buy now pay later delight shopper demi lovato reunites with exboyfriend joe jonas for a camp rock performance new pickleball kingdom indoor facility coming soon to south brunswick nj nasfm honor fire chief todd parker for excellence in emergency response virginia obamacare premium set to jump sharply scc report virginia obamacare premium set to jump sharply scc report going online in russia can be frustrating complicated and even dangerous chatgpt get smarter with gpt top feature youll love if youd invested in interactive broker ibkr stock year ago here how much youd have today graphene market to grow at cagr by lynch carpenter investigates claim in pandora data breach neptune frozen treat opening in malaga cove plaza top tech news today taylor swift ai scandal iphone price hike tesla dojo shutdown metas talent war siri upgrade stock market party to be ruined by sudden economic slowdown say stifel cnbc bitcoin hovers around whats next lockheed martin could lose major deal after u hit switzerland with tariff top meme coin to invest for maximum gain in xrp news today analyst react to m xrp transfer to coinbase is it routine move ford hit the pedal on ev production with billion overhaul of kentucky plant fed bowman say job data solidifies case for rate cut u consumer to bear brunt of tariff hit goldman economist say u consumer to bear brunt of tariff hit goldman economist say roll royce share price is on a sizzling hot streak but beware of the danger top laptop brand in the usa in ibd live qa stock list for monday aug why good leader ask great question manufacturer urged to act with technologyfirst plan a disruption intensifies say infotech research group in new report philippe cousteau earthecho international announces diverse youth advocate join the blue carbon ambassador program play the perfect game this summer with flightpath golf tee what to expect from tuesday report on inflation investopedia evercommerce to present at third quarter investor conference evercommerce to present at third quarter investor conference comedic multihyphenate matt mathews sign with wme hurry and grab a macbook air for just at back market air india suspends flight to washington dc amid ongoing tension this is what whale are betting on lockheed martin smart money is betting big in ups option aircraft gasket market is forecasted to reach u million in say stratview research best cryptos to invest in that could deliver massive longterm gain blockdag litecoin stellar toncoin solar stock outshine oil and gas benchmark after climate credit cut aerospace pilot control market is forecasted to reach u billion in say stratview research xml announces launch of new xml woman website relaunch of xml woman division aircraft gasket market is forecasted to reach u million in say stratview research aerospace pilot control market is forecasted to reach u billion in say stratview research diversified foodservice supply appoints lev peker a chief executive officer body armor plate market is forecasted to reach u million in say stratview research body armor plate market is forecasted to reach u million in say stratview research best typec hub for macbook in venchi bringing beloved italian chocolate to century city tegna stock rocket on nexstar merger talk fcc appears set to ease local tv rule in order to smooth deal path deadline btc miner stance on cloud mining in if youd invested in zoetis zts stock year ago here how much youd have today fda may not renew pfizers authorization for covid shot in younger child moderna prepares to boost supply ja mining offer fixed usd yield with crypto deposit alphabet inc nasdaqgoog share acquired by gouws capital llc kohmann bosshard financial service llc purchase share of costco wholesale corporation nasdaqcost gouws capital llc boost stock holding in thermo fisher scientific inc nysetmo ana alsharif became the face of the war in gaza for million then israel killed him home condition made rescue effort difficult pfd protect partner with sentin to pioneer aipowered industrial inspection gear up for the new school year with off therabody best seller wall street biggest bear is expecting a stockmarket pullback this year trump say gold import will not face u tariff bloombergcom irobots future isnt looking up russian seeking asylum in u facing new obstacle say usbased activist morse corp awarded m u army ota for uas development morse corp awarded m u army ota for uas development aem united state ag tractor and combine report july aem united state ag tractor and combine report july scientist say they cruised the ocean in a deepsea submersible and came across an undiscovered ecosystem macys is selling a cozy puma hoodie for thats super soft winnebago industry announces strategic leadership change to drive future growth how to stretch backtoschool budget at salvage store openai and google vet periodic lab seek million investment aol dialup got america online in the s it ending next month trump say gold will not be tariffed cnbc how turmeric water fight inflammation aussivo reveals vision to build blockchainverified cloud trust layer a new standard for transparency in the digital world the viral govee floor lamp is still at prime day price silver at doublevibration support ahead of gann time cluster unit building at rockaway reach completion in brownsville fill up your shopping cart with today best target deal including over off outdoor gear headphone kitchen good and more communication service stock whale activity in today session playhouse under construction for popup park in virginia beach aol is finally disconnecting it dialup internet year on barack obamas official portrait moved to nonpublic area in white house western union ceo ice raid tariff no threat to instamex deal video show moment of explosion at u steel plant in pennsylvania trump gross white house makeover reach petty new low crowded north county field vying to replace brian williams in missouri senate how the world most valuable company got caught in the middle of trump spat with china commentary vance right to call out warped partisan representation why is trump attacking intels ceo demi lovato reunites with exboyfriend joe jonas for a camp rock performance mamdani launch tour of new york city with a message linking cuomo to trump nasfm honor fire chief todd parker for excellence in emergency response what to know about the putintrump summit in alaska what to know about the putintrump summit in alaska what to know about the putintrump summit in alaska what to know about the putintrump summit in alaska what to know about the putintrump summit in alaska what to know about the putintrump summit in alaska what to know about the putintrump summit in alaska what to know about the putintrump summit in alaska what to know about the putintrump summit in alaska what to know about the putintrump summit in alaska what to know about the putintrump summit in alaska trump lash out at crockett renews call for cognitive test lynch carpenter investigates claim in pandora data breach president mayor police chief social worker trump is a man of many hat trump say he and putin will discus ukraine future border israel plan to widen coming offensive beyond gaza city into last area not under it control israel plan to widen coming offensive beyond gaza city into last area not under it control israel plan to widen coming offensive beyond gaza city into last area not under it control u senate candidate set the tone for kentucky republican during statewide dinner what to know about the putintrump summit in alaska what to know about the putintrump summit in alaska photo of trump emergency declaration a the national guard is activated in washington photo of trump emergency declaration a the national guard is activated in washington trump mistake alaska for russia and yes there snark surveillance case test chapter a standing and unfairness standard clinton county sheriff office arrest individual involved with multicounty sex trafficking operation a judge brutal rebuke of trump epstein gambit cnn politics what to know about the putintrump summit in alaska local view give ailing minnesota kid a fighting chance philippe cousteau earthecho international announces diverse youth advocate join the blue carbon ambassador program california couple file lawsuit against klp in conway over shortterm rental rule rubio hit medium over gaza narrative say some ally making it harder to achieve peace multiple men have impersonated ice agent to kidnap and assault woman multiple men have impersonated ice agent to kidnap and assault woman doe using ai dumb you down vote to save a historic building from demolition trump promised to fix inflation he doing this instead south carolina lawmaker awaiting trial on child sex abuse material charge resigns from office nvidia to share of revenue from chip sale to china with u government trump say gold wont face tariff trump say gold wont face tariff trump say gold wont face tariff trump say gold wont face tariff trump say gold wont face tariff ana alsharif became the face of the war in gaza for million then israel killed him trump say he placing washington police under federal control and deploying the national guard trump declares crime in dc will be cleaned up real quick during his announcement that he taking over the city police department chinese ship collide on chase in disputed water ap news summary at pm edt trump say he placing washington police under federal control and activating the national guard houlihan financial resource group ltd sell share of international business machine corporation nyseibm meridian wealth management llc ha million position in broadcom inc nasdaqavgo motley fool asset management llc raise stock holding in broadcom inc nasdaqavgo colorado congressional democrat visit aurora ice facility meet with detained activist jeanette vizguerra colorado congressional democrat visit aurora ice facility meet with detained activist jeanette vizguerra colorado congressional democrat visit aurora ice facility meet with detained activist jeanette vizguerra colorado congressional democrat visit aurora ice facility meet with detained activist jeanette vizguerra russian seeking asylum in u facing new obstacle say usbased activist trump warns protester amid dc crime crack down you spit and we hit trump suggests he can send national guard troop into other u city beyond dc scientist say they cruised the ocean in a deepsea submersible and came across an undiscovered ecosystem dillon associate inc increase holding in salesforce inc nysecrm colombian presidential candidate miguel uribe turbay dy two month after assassination attempt live now texas hold legislative session after democrat leave state jeffries hammer trump for seizing control of dc police no basis jeffries hammer trump for seizing control of dc police no basis jeffries hammer trump for seizing control of dc police no basis jeffries hammer trump for seizing control of dc police no basis jeffries hammer trump for seizing control of dc police no basis maga v miga laura loomer block marjorie taylor greene on social medium escalating a feud over america israel policy trump looking at reclassifying marijuana a less dangerous drug trump looking at reclassifying marijuana a less dangerous drug danish zoo asks for donated pet to feed to carnivore trump looking at reclassifying marijuana a less dangerous drug trump looking at reclassifying marijuana a less dangerous drug trump looking at reclassifying marijuana a less dangerous drug el paso police man dy when he hit tree get ejected trump team get blasted by judge for clumsy epstein diversion raw politics is now driving law enforcement at the department of justice civil right division nancy pelosi get all high and mighty lying about trump j and the ntl guard and oh hell no old woman this candidate for california governor ha a potential conflict of interest in her own home ishares core sp etf nysearcaivv share sold by aveo capital partner llc maine state library not victimized a suggested letter maine state library not victimized a suggested letter here are the country that have recognized a palestinian state johnson fistel investigates fairness of proposed sale of meridianlink currently the most dangerous place maia wealth llc purchase share of ishares core u aggregate bond etf nysearcaagg my wife of year just died my grief is so overwhelming i can barely cope
'@
